# Generate Report for Archive
#
# 1. Status text changed from "Ready for handoff" to "In Translation" on all
#    sheets that carry that status cell (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. The "Status" columns narrow (their header/value no longer needs as much
#    room once the text is shorter): Overview columns E & F, and column C on
#    both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the Status columns to match the new, shorter text ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
